$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- D8: add Fix Date (copy date formatting/style from an existing date cell) ---
$ws.Range("D8").Value = 44803
$ws.Range("A8").Copy()
$ws.Range("D8").PasteSpecial(-4122)

# --- Row 9: new error-log entry ---
$ws.Range("A9").Value = 44804
$ws.Range("A8").Copy()
$ws.Range("A9").PasteSpecial(-4122)

$ws.Range("B9").Value = "Charging/discharging final height needs to be more accurate. Bin isn't fully charging and discharging durring cycles."
$ws.Range("B2").Copy()
$ws.Range("B9").PasteSpecial(-4122)

$ws.Range("E9").Value = "FF.m"

$ws.Rows.Item(9).RowHeight = 31.5

# --- Row 10: new error-log entry ---
$ws.Range("A10").Value = 44804
$ws.Range("A8").Copy()
$ws.Range("A10").PasteSpecial(-4122)

$ws.Range("B10").Value = "Oscilations in temp profile are causing large temperature droops at the start of discharge"
$ws.Range("B2").Copy()
$ws.Range("B10").PasteSpecial(-4122)

$ws.Range("E10").Value = "FF.m"

$ws.Rows.Item(10).RowHeight = 31.5

# --- Update selection to match the author's final cursor position ---
$ws.Range("F13").Select()
